$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New row of data (row 4), mirroring the existing rows' structure.
$row = 4

$ws.Cells.Item($row, 1).Value = 7
$ws.Cells.Item($row, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item($row, 3).Value = "Ñuble"

# Date column (D) - keep the same date style/number format used in D2/D3.
# (set NumberFormat first so the engine reuses the existing style instead of
# minting a transient "general date" style before we override it)
$ws.Cells.Item($row, 4).NumberFormat = $ws.Cells.Item(3, 4).NumberFormat
$ws.Cells.Item($row, 4).Value = Get-Date -Year 2022 -Month 3 -Day 8 -Hour 0 -Minute 0 -Second 0 -Millisecond 0

$ws.Cells.Item($row, 5).Value = 16
$ws.Cells.Item($row, 6).Value = 100112043
$ws.Cells.Item($row, 7).Value = "Pepino dulce"
$ws.Cells.Item($row, 8).Value = "Cultivar IV Región"
$ws.Cells.Item($row, 9).Value = "Primera"
$ws.Cells.Item($row, 10).Value = 60
$ws.Cells.Item($row, 11).Value = 15000
$ws.Cells.Item($row, 12).Value = 16000
$ws.Cells.Item($row, 13).Value = 15500
$ws.Cells.Item($row, 14).Value = "`$/bandeja 18 kilos"
$ws.Cells.Item($row, 15).Value = "Provincia de Limarí"
$ws.Cells.Item($row, 16).Value = 861
$ws.Cells.Item($row, 17).Value = 18
$ws.Cells.Item($row, 18).Value = "Hortaliza"
